# NYPD CompStat 121st Precinct - weekly refresh to the following reporting
# week: "Volume 30 Number 25" -> "Volume 30 Number 26", and
# "6/19/2023 .. 6/25/2023" -> "6/26/2023 .. 7/2/2023", plus updated crime
# figures for the week.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Header / report period text (rich-text shared strings; plain .Value write
# is fine since every run in both strings shares the same font/format).
# ---------------------------------------------------------------------------
$ws.Range("A8").Value = "Volume 30   Number  26"
$ws.Range("C9").Value = "Report Covering the Week  6/26/2023  Through  7/2/2023"

# ---------------------------------------------------------------------------
# Helpers for cells that change *type* (number <-> "n/a" style text) so the
# cell style index matches the target exactly (the existing well-known
# "text" cells / "numeric" cells in the same table are reused as format
# donors via copy -> paste-formats, then the value is written).
# ---------------------------------------------------------------------------
function Set-TextZero([string]$addr) {
    $ws.Range("D14").Copy()
    $ws.Range($addr).PasteSpecial(-4122)   # xlPasteFormats
    $ws.Range("D14").Copy()
    $ws.Range($addr).PasteSpecial(-4163)   # xlPasteValues
}

function Set-TextNA([string]$addr) {
    $ws.Range("E14").Copy()
    $ws.Range($addr).PasteSpecial(-4122)   # xlPasteFormats
    $ws.Range("E14").Copy()
    $ws.Range($addr).PasteSpecial(-4163)   # xlPasteValues
}

function Set-NumCount([string]$addr, $value) {
    $ws.Range("F14").Copy()
    $ws.Range($addr).PasteSpecial(-4122)   # xlPasteFormats
    $ws.Range($addr).Value = $value
}

function Set-NumPct([string]$addr, $value) {
    $ws.Range("L14").Copy()
    $ws.Range($addr).PasteSpecial(-4122)   # xlPasteFormats
    $ws.Range($addr).Value = $value
}

# ---------------------------------------------------------------------------
# Row 14 - Murder
# ---------------------------------------------------------------------------
Set-TextZero "C14"

# ---------------------------------------------------------------------------
# Row 15 - Rape
# ---------------------------------------------------------------------------
Set-TextZero "G15"
Set-TextNA   "H15"

# ---------------------------------------------------------------------------
# Row 16 - Robbery
# ---------------------------------------------------------------------------
Set-TextZero "C16"
$ws.Range("F16").Value = 7
$ws.Range("G16").Value = 2
$ws.Range("H16").Value = 250
$ws.Range("L16").Value = 27.586206896551

# ---------------------------------------------------------------------------
# Row 17 - Fel. Assault
# ---------------------------------------------------------------------------
$ws.Range("C17").Value = 7
$ws.Range("D17").Value = 8
$ws.Range("E17").Value = -12.5
$ws.Range("F17").Value = 21
$ws.Range("G17").Value = 22
$ws.Range("H17").Value = -4.545454545454
$ws.Range("I17").Value = 122
$ws.Range("J17").Value = 128
$ws.Range("K17").Value = -4.6875
$ws.Range("L17").Value = 60.526315789473

# ---------------------------------------------------------------------------
# Row 18 - Burglary
# ---------------------------------------------------------------------------
Set-TextZero "D18"
Set-TextNA   "E18"
$ws.Range("F18").Value = 7
$ws.Range("G18").Value = 2
$ws.Range("H18").Value = 250
$ws.Range("I18").Value = 26
$ws.Range("K18").Value = -33.333333333333
$ws.Range("L18").Value = -44.680851063829

# ---------------------------------------------------------------------------
# Row 19 - Gr. Larceny
# ---------------------------------------------------------------------------
$ws.Range("C19").Value = 8
$ws.Range("D19").Value = 2
$ws.Range("E19").Value = 300
$ws.Range("F19").Value = 38
$ws.Range("G19").Value = 22
$ws.Range("H19").Value = 72.727272727272
$ws.Range("I19").Value = 196
$ws.Range("J19").Value = 227
$ws.Range("K19").Value = -13.656387665198
$ws.Range("L19").Value = 37.062937062937

# ---------------------------------------------------------------------------
# Row 20 - G.L.A.
# ---------------------------------------------------------------------------
$ws.Range("C20").Value = 2
$ws.Range("D20").Value = 2
$ws.Range("E20").Value = 0
$ws.Range("G20").Value = 6
$ws.Range("H20").Value = 66.666666666666
$ws.Range("I20").Value = 56
$ws.Range("J20").Value = 42
$ws.Range("K20").Value = 33.333333333333
$ws.Range("L20").Value = 211.111111111111

# ---------------------------------------------------------------------------
# Row 21 - TOTAL
# ---------------------------------------------------------------------------
$ws.Range("C21").Value = 19
$ws.Range("D21").Value = 12
$ws.Range("E21").Value = 58.333333333333
$ws.Range("G21").Value = 54
$ws.Range("H21").Value = 57.407407407407
$ws.Range("I21").Value = 443
$ws.Range("J21").Value = 470
$ws.Range("K21").Value = -5.744680851063
$ws.Range("L21").Value = 36.307692307692

# ---------------------------------------------------------------------------
# Row 23 - Housing
# ---------------------------------------------------------------------------
$ws.Range("L23").Value = 11.111111111111

# ---------------------------------------------------------------------------
# Row 24 - Petit Larceny
# ---------------------------------------------------------------------------
$ws.Range("C24").Value = 32
$ws.Range("D24").Value = 28
$ws.Range("E24").Value = 14.285714285714
$ws.Range("F24").Value = 133
$ws.Range("G24").Value = 149
$ws.Range("H24").Value = -10.738255033557
$ws.Range("I24").Value = 687
$ws.Range("J24").Value = 674
$ws.Range("K24").Value = 1.928783382789
$ws.Range("L24").Value = 46.794871794871

# ---------------------------------------------------------------------------
# Row 25 - Misd. Assault
# ---------------------------------------------------------------------------
$ws.Range("C25").Value = 8
$ws.Range("D25").Value = 19
$ws.Range("E25").Value = -57.894736842105
$ws.Range("G25").Value = 41
$ws.Range("H25").Value = -7.317073170731
$ws.Range("I25").Value = 276
$ws.Range("J25").Value = 249
$ws.Range("K25").Value = 10.843373493975
$ws.Range("L25").Value = 46.031746031746

# ---------------------------------------------------------------------------
# Row 26 - UCR Rape*
# ---------------------------------------------------------------------------
$ws.Range("G26").Value = 1
$ws.Range("H26").Value = 0

# ---------------------------------------------------------------------------
# Row 27 - Other Sex Crimes
# ---------------------------------------------------------------------------
Set-NumCount "D27" 1
Set-NumPct   "E27" -100
$ws.Range("F27").Value = 3
$ws.Range("G27").Value = 9
$ws.Range("H27").Value = -66.666666666666
$ws.Range("J27").Value = 27
$ws.Range("K27").Value = 33.333333333333
$ws.Range("L27").Value = 80

# ---------------------------------------------------------------------------
# Row 30 - Hate Crimes
# ---------------------------------------------------------------------------
Set-NumCount "D30" 1
Set-NumPct   "E30" -100
Set-NumCount "G30" 1
Set-NumPct   "H30" 100
$ws.Range("J30").Value = 6
$ws.Range("K30").Value = -50
